# Weekly data refresh: "Add data for 2022-07-17"
# Updates year-to-date violent-crime counts (mostly the 2022 column, plus a
# few small corrections to earlier years) across the citywide, by-neighborhood
# summary, and per-neighborhood detail sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 63
$ws.Range("E3").Value = 82
$ws.Range("F3").Value = 77
$ws.Range("I3").Value = 103
$ws.Range("B6").Value = 210
$ws.Range("D6").Value = 242
$ws.Range("E6").Value = 235
$ws.Range("F6").Value = 297
$ws.Range("G6").Value = 266
$ws.Range("H6").Value = 234
$ws.Range("I6").Value = 297
$ws.Range("B7").Value = 283
$ws.Range("D7").Value = 376
$ws.Range("E7").Value = 363
$ws.Range("F7").Value = 426
$ws.Range("G7").Value = 392
$ws.Range("H7").Value = 364
$ws.Range("I7").Value = 479
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I2").Value = 2
$ws.Range("I6").Value = 11
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("D4").Value = 5
$ws.Range("D5").Value = 7
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("E3").Value = 3
$ws.Range("E6").Value = 5
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B5").Value = 23
$ws.Range("H5").Value = 19
$ws.Range("I5").Value = 12
$ws.Range("B6").Value = 23
$ws.Range("H6").Value = 32
$ws.Range("I6").Value = 27
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I3").Value = 2
$ws.Range("I5").Value = 9
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I5").Value = 11
$ws.Range("E8").Value = 23
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = 8
$ws.Range("D19").Value = 8
$ws.Range("F20").Value = 8
$ws.Range("B27").Value = 23
$ws.Range("H27").Value = 32
$ws.Range("I27").Value = 27
$ws.Range("H47").Value = 2
$ws.Range("G52").Value = 50
$ws.Range("I52").Value = 73
$ws.Range("G69").Value = 11
$ws.Range("G73").Value = 7
$ws.Range("F75").Value = 12
$ws.Range("I81").Value = 9
$ws.Range("E85").Value = 5
$ws.Range("D87").Value = 7
$ws.Range("B97").Value = 283
$ws.Range("D97").Value = 376
$ws.Range("E97").Value = 363
$ws.Range("F97").Value = 426
$ws.Range("G97").Value = 392
$ws.Range("H97").Value = 364
$ws.Range("I97").Value = 479
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 8
$ws.Range("G6").Value = 34
$ws.Range("G7").Value = 50
$ws.Range("I7").Value = 73
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("F5").Value = 10
$ws.Range("F6").Value = 12
$ws = $wb.Worksheets.Item("River North")
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 7
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("F3").Value = 3
$ws.Range("H5").Value = 5
$ws.Range("F6").Value = 7
$ws.Range("H6").Value = 8
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 2
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E5").Value = 15
$ws.Range("E6").Value = 23
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("F4").Value = 8
$ws.Range("F5").Value = 11
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("F3").Value = 1
$ws.Range("F7").Value = 8
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("D5").Value = 6
$ws.Range("D6").Value = 8
